# Commit: "edited request command messages"
#
# On the Undo/Redo activity-diagram slide, two shapes mention the (old)
# "order book" / "orderBookStateList" terminology from the request-command
# flow; rename them to "address book" / "addressBookStateList".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "TextBox 47": "[command commits order book]" -> "[command commits address book]"
$shpGuard = $s.Shapes.Item("TextBox 47")
$trGuard = $shpGuard.TextFrame.TextRange
$fullGuard = $trGuard.Text
$oldGuardTail = "command commits order book]"
$newGuardTail = "command commits address book]"
$tailStart = $fullGuard.IndexOf($oldGuardTail) + 1
if ($tailStart -gt 0) {
    $tailRange = $trGuard.Characters($tailStart, $oldGuardTail.Length)
    $tailRange.Text = $newGuardTail
}

# --- "Rounded Rectangle 50": purge-state note mentions "order book" / "orderBookStateList"
$shpNote = $s.Shapes.Item("Rounded Rectangle 50")
$trNote = $shpNote.TextFrame.TextRange

# First collapse "...save order book to " -> "...save address book to " (this
# also merges the run boundary that used to split before "order book").
$fullNote = $trNote.Text
$oldHead = "Purge redundant states and then save order book to "
$newHead = "Purge redundant states and then save address book to "
$headStart = $fullNote.IndexOf($oldHead) + 1
if ($headStart -gt 0) {
    $headRange = $trNote.Characters($headStart, $oldHead.Length)
    $headRange.Text = $newHead
}

# Then rename the "orderBookStateList" identifier run to "addressBookStateList".
$trNote2 = $shpNote.TextFrame.TextRange
$fullNote2 = $trNote2.Text
$oldId = "orderBookStateList"
$newId = "addressBookStateList"
$idStart = $fullNote2.IndexOf($oldId) + 1
if ($idStart -gt 0) {
    $idRange = $trNote2.Characters($idStart, $oldId.Length)
    $idRange.Text = $newId
}
